# ---------------------------------------------------------------------------
# Adds a new "2022-Q1" sheet (fund holders detail) positioned right before
# the existing "总计" (totals) sheet, and inserts a corresponding summary
# row at the top of the "总计" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q1" worksheet, placed right before "总计" ----
$totalSheet = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($totalSheet)
$ws.Name = "2022-Q1"

# Copy the header-row / index-column formatting (bold / centered / boxed,
# the style used by the other quarterly sheets) from an existing sheet so
# the new tab matches the established look.
$styleSource = $wb.Worksheets.Item("2021-Q4")
$styleSource.Range("B1:H1").Copy($ws.Range("B1:H1"))
$styleSource.Range("A2:A29").Copy($ws.Range("A2:A29"))

# --- 2. Header row -----------------------------------------------------
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

# --- 3. Data rows (28 funds) --------------------------------------------
# Columns B-G are stored as text (to preserve leading zeros in fund codes
# and fixed decimal formatting); column A (row index) and H (rank) are
# numeric.
$dataText = @"
0	166019	中欧价值智选回报混合A	156.17	94.14	5.45	8.5113	4
1	013220	中欧新兴价值一年持有混合A	63.33	94.47	6.84	4.3318	3
2	004235	中欧价值智选回报混合C	36.40	94.14	5.45	1.9838	4
3	004848	中欧睿泓定期开放灵活配置混合	23.30	59.08	6.53	1.5215	2
4	013221	中欧新兴价值一年持有混合C	16.75	94.47	6.84	1.1457	3
5	001887	中欧价值智选回报混合E	20.77	94.14	5.45	1.1320	4
6	005821	万家新机遇龙头企业灵活配置混合	23.23	56.20	3.54	0.8223	2
7	270028	广发制造业精选混合A	26.52	93.12	2.63	0.6975	9
8	013960	万家新机遇成长一年持有期混合A	13.29	49.51	3.69	0.4904	1
9	000327	南方潜力新蓝筹混合	6.58	93.61	6.95	0.4573	2
10	003516	国泰融安多策略灵活配置混合	11.18	71.30	2.76	0.3086	4
11	168501	北信瑞丰产业升级多策略混合	4.42	94.11	4.25	0.1878	6
12	005412	金信民长灵活配置混合A	2.37	70.53	6.33	0.1500	2
13	013961	万家新机遇成长一年持有期混合C	3.13	49.51	3.69	0.1155	1
14	020023	国泰事件驱动策略混合	3.34	74.02	2.95	0.0985	3
15	010023	广发制造业精选混合C	3.57	93.12	2.63	0.0939	9
16	005413	金信民长灵活配置混合C	1.30	70.53	6.33	0.0823	2
17	005894	华夏优势精选股票	1.92	91.06	3.81	0.0732	8
18	001056	北信瑞丰健康生活主题灵活配置混合	1.64	86.03	3.82	0.0626	6
19	008180	同泰慧利混合A	1.40	92.38	4.03	0.0564	8
20	012445	华富新能源股票型发起式证券投资基金	1.35	85.55	3.64	0.0491	8
21	229002	泰达宏利逆向策略混合	1.63	92.33	1.53	0.0249	6
22	001017	泰达宏利改革动力量化策略灵活配置混合A	1.50	92.25	1.64	0.0246	8
23	000679	招商丰利灵活配置混合A	0.39	74.75	4.52	0.0176	3
24	005903	泰达宏利绩优增长灵活配置混合	0.40	92.75	3.69	0.0148	8
25	008181	同泰慧利混合C	0.19	92.38	4.03	0.0077	8
26	002416	招商丰利灵活配置混合C	0.02	74.75	4.52	0.0009	3
27	003550	泰达宏利改革动力量化策略灵活配置混合C	0.01	92.25	1.64	0.0002	8
"@

$lines = $dataText -split "`n"

$ws.Range("B2:G29").NumberFormat = "@"

$rowNum = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $cols = $line -split "`t"

    $ws.Cells.Item($rowNum, 1).Value = [int]$cols[0]
    $ws.Cells.Item($rowNum, 2).Value = $cols[1]
    $ws.Cells.Item($rowNum, 3).Value = $cols[2]
    $ws.Cells.Item($rowNum, 4).Value = $cols[3]
    $ws.Cells.Item($rowNum, 5).Value = $cols[4]
    $ws.Cells.Item($rowNum, 6).Value = $cols[5]
    $ws.Cells.Item($rowNum, 7).Value = $cols[6]
    $ws.Cells.Item($rowNum, 8).Value = [int]$cols[7]

    $rowNum++
}

# --- 4. Update the "总计" (totals) sheet: insert the 2022-Q1 summary ----
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert()

$srcA = $totalWs.Cells.Item(3, 1)
$dstA = $totalWs.Cells.Item(2, 1)
$srcA.Copy($dstA)
$dstA.Value = 0

$totalWs.Range("B2:D2").ClearFormats()
$totalWs.Cells.Item(2, 2).Value = "2022-Q1"
$totalWs.Cells.Item(2, 3).Value = 28
$totalWs.Cells.Item(2, 4).Value = 22.46

for ($r = 3; $r -le 7; $r++) {
    $totalWs.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "2022-Q1 sheet added and 总计 sheet updated."
